$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("G2").Value = [double]"20.32821866666667"
$ws.Range("H2").Value = [double]"60.984656"
$ws.Range("I2").Value = [double]"0.004181898474048532"
$ws.Range("J2").Value = [double]"0.004181898474048532"
$ws.Range("K2").Value = [double]"3"
$ws.Range("M2").Value = [double]"169.629438"
$ws.Range("N2").Value = [double]"508.888314"
$ws.Range("O2").Value = [double]"0.7428377317484701"
$ws.Range("P2").Value = [double]"0.7428377317484702"
$ws.Range("Q2").Value = [double]"3448.264307967776"
$ws.Range("R2").Value = [double]"31034.37877170998"
$ws.Range("S2").Value = [double]"0.0031064719768646"
$ws.Range("T2").Value = [double]"0.0031064719768646"
$ws.Range("E3").Value = [double]"3"
$ws.Range("G3").Value = [double]"20.32821866666667"
$ws.Range("H3").Value = [double]"60.984656"
$ws.Range("I3").Value = [double]"0.004181898474048532"
$ws.Range("J3").Value = [double]"0.004181898474048532"
$ws.Range("K3").Value = [double]"3"
$ws.Range("M3").Value = [double]"0.9848756666666668"
$ws.Range("N3").Value = [double]"2.954627"
$ws.Range("O3").Value = [double]"0.004312947180081616"
$ws.Range("P3").Value = [double]"0.004312947180081616"
$ws.Range("Q3").Value = [double]"20.02076791147912"
$ws.Range("R3").Value = [double]"180.186911203312"
$ws.Range("S3").Value = [double]"1.803630723103523E-05"
$ws.Range("T3").Value = [double]"1.803630723103523E-05"
$ws.Range("E4").Value = [double]"3"
$ws.Range("G4").Value = [double]"20.32821866666667"
$ws.Range("H4").Value = [double]"60.984656"
$ws.Range("I4").Value = [double]"0.004181898474048532"
$ws.Range("J4").Value = [double]"0.004181898474048532"
$ws.Range("K4").Value = [double]"3"
$ws.Range("M4").Value = [double]"54.620752"
$ws.Range("N4").Value = [double]"163.862256"
$ws.Range("O4").Value = [double]"0.2391940691454494"
$ws.Range("P4").Value = [double]"0.2391940691454494"
$ws.Range("Q4").Value = [double]"1110.342590393771"
$ws.Range("R4").Value = [double]"9993.083313543937"
$ws.Range("S4").Value = [double]"0.001000285312760814"
$ws.Range("T4").Value = [double]"0.001000285312760814"
$ws.Range("E5").Value = [double]"3"
$ws.Range("G5").Value = [double]"20.32821866666667"
$ws.Range("H5").Value = [double]"60.984656"
$ws.Range("I5").Value = [double]"0.004181898474048532"
$ws.Range("J5").Value = [double]"0.004181898474048532"
$ws.Range("K5").Value = [double]"3"
$ws.Range("M5").Value = [double]"3.118221666666667"
$ws.Range("N5").Value = [double]"9.354665000000001"
$ws.Range("O5").Value = [double]"0.01365525192599884"
$ws.Range("P5").Value = [double]"0.01365525192599884"
$ws.Range("Q5").Value = [double]"63.38789189113778"
$ws.Range("R5").Value = [double]"570.49102702024"
$ws.Range("S5").Value = [double]"5.710487719208283E-05"
$ws.Range("T5").Value = [double]"5.710487719208283E-05"
$ws.Range("E6").Value = [double]"3"
$ws.Range("G6").Value = [double]"4809.896321333334"
$ws.Range("H6").Value = [double]"14429.688964"
$ws.Range("I6").Value = [double]"0.9894865072215304"
$ws.Range("J6").Value = [double]"0.9894865072215304"
$ws.Range("K6").Value = [double]"3"
$ws.Range("M6").Value = [double]"169.629438"
$ws.Range("N6").Value = [double]"508.888314"
$ws.Range("O6").Value = [double]"0.7428377317484701"
$ws.Range("P6").Value = [double]"0.7428377317484702"
$ws.Range("Q6").Value = [double]"815900.0098260407"
$ws.Range("R6").Value = [double]"7343100.088434367"
$ws.Range("S6").Value = [double]"0.7350279126201579"
$ws.Range("T6").Value = [double]"0.735027912620158"
$ws.Range("E7").Value = [double]"3"
$ws.Range("G7").Value = [double]"4809.896321333334"
$ws.Range("H7").Value = [double]"14429.688964"
$ws.Range("I7").Value = [double]"0.9894865072215304"
$ws.Range("J7").Value = [double]"0.9894865072215304"
$ws.Range("K7").Value = [double]"3"
$ws.Range("M7").Value = [double]"0.9848756666666668"
$ws.Range("N7").Value = [double]"2.954627"
$ws.Range("O7").Value = [double]"0.004312947180081616"
$ws.Range("P7").Value = [double]"0.004312947180081616"
$ws.Range("Q7").Value = [double]"4737.149846070715"
$ws.Range("R7").Value = [double]"42634.34861463644"
$ws.Range("S7").Value = [double]"0.004267603041049907"
$ws.Range("T7").Value = [double]"0.004267603041049907"
$ws.Range("E8").Value = [double]"3"
$ws.Range("G8").Value = [double]"4809.896321333334"
$ws.Range("H8").Value = [double]"14429.688964"
$ws.Range("I8").Value = [double]"0.9894865072215304"
$ws.Range("J8").Value = [double]"0.9894865072215304"
$ws.Range("K8").Value = [double]"3"
$ws.Range("M8").Value = [double]"54.620752"
$ws.Range("N8").Value = [double]"163.862256"
$ws.Range("O8").Value = [double]"0.2391940691454494"
$ws.Range("P8").Value = [double]"0.2391940691454494"
$ws.Range("Q8").Value = [double]"262720.1541132603"
$ws.Range("R8").Value = [double]"2364481.387019343"
$ws.Range("S8").Value = [double]"0.236679304026836"
$ws.Range("T8").Value = [double]"0.236679304026836"
$ws.Range("E9").Value = [double]"3"
$ws.Range("G9").Value = [double]"4809.896321333334"
$ws.Range("H9").Value = [double]"14429.688964"
$ws.Range("I9").Value = [double]"0.9894865072215304"
$ws.Range("J9").Value = [double]"0.9894865072215304"
$ws.Range("K9").Value = [double]"3"
$ws.Range("M9").Value = [double]"3.118221666666667"
$ws.Range("N9").Value = [double]"9.354665000000001"
$ws.Range("O9").Value = [double]"0.01365525192599884"
$ws.Range("P9").Value = [double]"0.01365525192599884"
$ws.Range("Q9").Value = [double]"14998.3229236019"
$ws.Range("R9").Value = [double]"134984.9063124171"
$ws.Range("S9").Value = [double]"0.01351168753348667"
$ws.Range("T9").Value = [double]"0.01351168753348667"
$ws.Range("E10").Value = [double]"3"
$ws.Range("G10").Value = [double]"2.69506"
$ws.Range("H10").Value = [double]"8.085180000000001"
$ws.Range("I10").Value = [double]"0.000554424737665286"
$ws.Range("J10").Value = [double]"0.000554424737665286"
$ws.Range("K10").Value = [double]"3"
$ws.Range("M10").Value = [double]"169.629438"
$ws.Range("N10").Value = [double]"508.888314"
$ws.Range("O10").Value = [double]"0.7428377317484701"
$ws.Range("P10").Value = [double]"0.7428377317484702"
$ws.Range("Q10").Value = [double]"457.16151317628"
$ws.Range("R10").Value = [double]"4114.453618586521"
$ws.Range("S10").Value = [double]"0.0004118476145525217"
$ws.Range("T10").Value = [double]"0.0004118476145525217"
$ws.Range("E11").Value = [double]"3"
$ws.Range("G11").Value = [double]"2.69506"
$ws.Range("H11").Value = [double]"8.085180000000001"
$ws.Range("I11").Value = [double]"0.000554424737665286"
$ws.Range("J11").Value = [double]"0.000554424737665286"
$ws.Range("K11").Value = [double]"3"
$ws.Range("M11").Value = [double]"0.9848756666666668"
$ws.Range("N11").Value = [double]"2.954627"
$ws.Range("O11").Value = [double]"0.004312947180081616"
$ws.Range("P11").Value = [double]"0.004312947180081616"
$ws.Range("Q11").Value = [double]"2.654299014206667"
$ws.Range("R11").Value = [double]"23.88869112786001"
$ws.Range("S11").Value = [double]"2.391204608880985E-06"
$ws.Range("T11").Value = [double]"2.391204608880985E-06"
$ws.Range("E12").Value = [double]"3"
$ws.Range("G12").Value = [double]"2.69506"
$ws.Range("H12").Value = [double]"8.085180000000001"
$ws.Range("I12").Value = [double]"0.000554424737665286"
$ws.Range("J12").Value = [double]"0.000554424737665286"
$ws.Range("K12").Value = [double]"3"
$ws.Range("M12").Value = [double]"54.620752"
$ws.Range("N12").Value = [double]"163.862256"
$ws.Range("O12").Value = [double]"0.2391940691454494"
$ws.Range("P12").Value = [double]"0.2391940691454494"
$ws.Range("Q12").Value = [double]"147.20620388512"
$ws.Range("R12").Value = [double]"1324.85583496608"
$ws.Range("S12").Value = [double]"0.0001326151090370581"
$ws.Range("T12").Value = [double]"0.0001326151090370581"
$ws.Range("E13").Value = [double]"3"
$ws.Range("G13").Value = [double]"2.69506"
$ws.Range("H13").Value = [double]"8.085180000000001"
$ws.Range("I13").Value = [double]"0.000554424737665286"
$ws.Range("J13").Value = [double]"0.000554424737665286"
$ws.Range("K13").Value = [double]"3"
$ws.Range("M13").Value = [double]"3.118221666666667"
$ws.Range("N13").Value = [double]"9.354665000000001"
$ws.Range("O13").Value = [double]"0.01365525192599884"
$ws.Range("P13").Value = [double]"0.01365525192599884"
$ws.Range("Q13").Value = [double]"8.403794484966667"
$ws.Range("R13").Value = [double]"75.63415036470002"
$ws.Range("S13").Value = [double]"7.5708094668253E-06"
$ws.Range("T13").Value = [double]"7.5708094668253E-06"
$ws.Range("E14").Value = [double]"3"
$ws.Range("G14").Value = [double]"28.08283533333333"
$ws.Range("H14").Value = [double]"84.24850599999999"
$ws.Range("I14").Value = [double]"0.005777169566755752"
$ws.Range("J14").Value = [double]"0.005777169566755752"
$ws.Range("K14").Value = [double]"3"
$ws.Range("M14").Value = [double]"169.629438"
$ws.Range("N14").Value = [double]"508.888314"
$ws.Range("O14").Value = [double]"0.7428377317484701"
$ws.Range("P14").Value = [double]"0.7428377317484702"
$ws.Range("Q14").Value = [double]"4763.675575039875"
$ws.Range("R14").Value = [double]"42873.08017535887"
$ws.Range("S14").Value = [double]"0.004291499536895135"
$ws.Range("T14").Value = [double]"0.004291499536895135"
$ws.Range("E15").Value = [double]"3"
$ws.Range("G15").Value = [double]"28.08283533333333"
$ws.Range("H15").Value = [double]"84.24850599999999"
$ws.Range("I15").Value = [double]"0.005777169566755752"
$ws.Range("J15").Value = [double]"0.005777169566755752"
$ws.Range("K15").Value = [double]"3"
$ws.Range("M15").Value = [double]"0.9848756666666668"
$ws.Range("N15").Value = [double]"2.954627"
$ws.Range("O15").Value = [double]"0.004312947180081616"
$ws.Range("P15").Value = [double]"0.004312947180081616"
$ws.Range("Q15").Value = [double]"27.65810117080689"
$ws.Range("R15").Value = [double]"248.922910537262"
$ws.Range("S15").Value = [double]"2.491662719179255E-05"
$ws.Range("T15").Value = [double]"2.491662719179255E-05"
$ws.Range("E16").Value = [double]"3"
$ws.Range("G16").Value = [double]"28.08283533333333"
$ws.Range("H16").Value = [double]"84.24850599999999"
$ws.Range("I16").Value = [double]"0.005777169566755752"
$ws.Range("J16").Value = [double]"0.005777169566755752"
$ws.Range("K16").Value = [double]"3"
$ws.Range("M16").Value = [double]"54.620752"
$ws.Range("N16").Value = [double]"163.862256"
$ws.Range("O16").Value = [double]"0.2391940691454494"
$ws.Range("P16").Value = [double]"0.2391940691454494"
$ws.Range("Q16").Value = [double]"1533.905584198837"
$ws.Range("R16").Value = [double]"13805.15025778954"
$ws.Range("S16").Value = [double]"0.001381864696815561"
$ws.Range("T16").Value = [double]"0.001381864696815561"
$ws.Range("E17").Value = [double]"3"
$ws.Range("G17").Value = [double]"28.08283533333333"
$ws.Range("H17").Value = [double]"84.24850599999999"
$ws.Range("I17").Value = [double]"0.005777169566755752"
$ws.Range("J17").Value = [double]"0.005777169566755752"
$ws.Range("K17").Value = [double]"3"
$ws.Range("M17").Value = [double]"3.118221666666667"
$ws.Range("N17").Value = [double]"9.354665000000001"
$ws.Range("O17").Value = [double]"0.01365525192599884"
$ws.Range("P17").Value = [double]"0.01365525192599884"
$ws.Range("Q17").Value = [double]"87.56850559783223"
$ws.Range("R17").Value = [double]"788.11655038049"
$ws.Range("S17").Value = [double]"7.888870585326338E-05"
$ws.Range("T17").Value = [double]"7.888870585326338E-05"
